$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_info")

# Shift the whole table from A1:C8 down-right to B2:D9
$ws.Range("A1").EntireColumn.Insert()
$ws.Range("A1").EntireRow.Insert()

# Fix typo in the "data" row description (now at D4)
$ws.Range("D4").Value = "Work date for the timecard entry (YYYY-MM-DD)"

# Header row formatting (B2:D2)
$hdr = $ws.Range("B2:D2")
$hdr.Font.Bold = $true
$hdr.Font.Size = 14
$hdr.Font.Name = "Aptos Narrow"
$hdr.Interior.ThemeColor = 7
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true
$ws.Rows.Item(2).RowHeight = 18

# Body formatting (B3:D9)
$body = $ws.Range("B3:D9")
$body.Borders.LineStyle = 1
$body.VerticalAlignment = -4108
$body.WrapText = $true

# Sheet view tweaks
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("G15").Select()

Write-Output "done"
